$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the trailing footer/metadata rows (1001-1005). This also
#    shrinks the worksheet dimension from A1:D1005 down to A1:D999.
$ws.Rows("1001:1005").Delete()

# 2. Title-case the Spanish connector words ("de", "del", "el", "los",
#    "la", "las") that appear inside the state/municipality names.
#    Every occurrence in the sheet is surrounded by plain spaces, so a
#    simple substring replace is equivalent to a whole-word replace.
[void]$ws.Cells.Replace(" de ", " De ")
[void]$ws.Cells.Replace(" del ", " Del ")
[void]$ws.Cells.Replace(" el ", " El ")
[void]$ws.Cells.Replace(" los ", " Los ")
[void]$ws.Cells.Replace(" la ", " La ")
[void]$ws.Cells.Replace(" las ", " Las ")

# 3. Rename the header row to the new short English column names.
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 4. Fix the tiny floating point rounding differences for the
#    "6 matriculas" percentage (0.0009837678307919331 -> ...333), which
#    affects every row whose "Numero de Matriculas" equals 6.
$fixedCells = @("D4", "D31", "D51", "D74", "D82", "D88", "D92", "D125", "D186", "D241", "D253", "D261", "D280", "D300", "D388", "D389", "D499", "D512", "D521", "D562", "D579", "D709", "D717", "D738", "D782", "D822", "D851", "D979")
foreach ($addr in $fixedCells) {
    $ws.Range($addr).Value = 0.0009837678307919333
}
